# Reorders Field Name / Field Modifier / Field Type values within each
# class-group on the "classFields" sheet (standard relationship + MSM measure work).
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("classFields")

$ws.Cells.Item(2, 2).Value = "name"
$ws.Cells.Item(2, 3).Value = "private"
$ws.Cells.Item(2, 4).Value = "java.lang.String"
$ws.Cells.Item(3, 2).Value = "`$VALUES"
$ws.Cells.Item(3, 3).Value = "private"
$ws.Cells.Item(3, 4).Value = "org.andante.enums.KafkaConsumerGroup[]"
$ws.Cells.Item(5, 2).Value = "PRODUCT_ORDER_GROUP"
$ws.Cells.Item(5, 3).Value = "public"
$ws.Cells.Item(5, 4).Value = "org.andante.enums.KafkaConsumerGroup"
$ws.Cells.Item(6, 2).Value = "ACTIVITY_ORDER_GROUP"
$ws.Cells.Item(6, 3).Value = "public"
$ws.Cells.Item(6, 4).Value = "org.andante.enums.KafkaConsumerGroup"
$ws.Cells.Item(7, 2).Value = "EUREKA_PORT"
$ws.Cells.Item(7, 4).Value = "java.lang.Integer"
$ws.Cells.Item(8, 2).Value = "kafkaContainer"
$ws.Cells.Item(8, 4).Value = "org.testcontainers.containers.KafkaContainer"
$ws.Cells.Item(9, 2).Value = "eurekaContainer"
$ws.Cells.Item(9, 4).Value = "org.testcontainers.containers.GenericContainer"
$ws.Cells.Item(11, 2).Value = "`$VALUES"
$ws.Cells.Item(11, 3).Value = "private"
$ws.Cells.Item(11, 4).Value = "org.andante.enums.OperationStatus[]"
$ws.Cells.Item(12, 2).Value = "NOT_FOUND"
$ws.Cells.Item(12, 3).Value = "public"
$ws.Cells.Item(12, 4).Value = "org.andante.enums.OperationStatus"
$ws.Cells.Item(15, 2).Value = "PRODUCT_SUBWOOFERS_VARIANT_INTERNAL_TOPIC"
$ws.Cells.Item(16, 2).Value = "PRODUCT_COMMENT_INTERNAL_TOPIC"
$ws.Cells.Item(17, 2).Value = "ORDER_INTERNAL_TOPIC"
$ws.Cells.Item(18, 2).Value = "PRODUCT_AMPLIFIER_VARIANT_INTERNAL_TOPIC"
$ws.Cells.Item(18, 3).Value = "public"
$ws.Cells.Item(18, 4).Value = "org.andante.enums.KafkaTopic"
$ws.Cells.Item(19, 2).Value = "PRODUCT_SPEAKERS_INTERNAL_TOPIC"
$ws.Cells.Item(20, 2).Value = "topicName"
$ws.Cells.Item(20, 3).Value = "private"
$ws.Cells.Item(20, 4).Value = "java.lang.String"
$ws.Cells.Item(21, 2).Value = "`$VALUES"
$ws.Cells.Item(21, 3).Value = "private"
$ws.Cells.Item(21, 4).Value = "org.andante.enums.KafkaTopic[]"
$ws.Cells.Item(22, 2).Value = "PRODUCT_HEADPHONES_INTERNAL_TOPIC"
$ws.Cells.Item(23, 2).Value = "PRODUCT_ORDER_INTERNAL_TOPIC"
$ws.Cells.Item(24, 2).Value = "ORDER_ENTRY_INTERNAL_TOPIC"
$ws.Cells.Item(25, 2).Value = "PRODUCT_HEADPHONES_VARIANT_INTERNAL_TOPIC"
$ws.Cells.Item(26, 2).Value = "PRODUCT_MICROPHONE_INTERNAL_TOPIC"
$ws.Cells.Item(27, 2).Value = "PRODUCT_SPEAKERS_VARIANT_INTERNAL_TOPIC"
$ws.Cells.Item(28, 2).Value = "PRODUCT_GRAMOPHONE_INTERNAL_TOPIC"
$ws.Cells.Item(29, 2).Value = "PRODUCT_GRAMOPHONE_VARIANT_INTERNAL_TOPIC"
$ws.Cells.Item(30, 2).Value = "PRODUCT_AMPLIFIER_INTERNAL_TOPIC"
$ws.Cells.Item(31, 2).Value = "PRODUCT_SUBWOOFERS_INTERNAL_TOPIC"
$ws.Cells.Item(32, 2).Value = "PRODUCT_PRODUCER_INTERNAL_TOPIC"
$ws.Cells.Item(33, 2).Value = "PRODUCT_MICROPHONE_VARIANT_INTERNAL_TOPIC"
$ws.Cells.Item(33, 3).Value = "public"
$ws.Cells.Item(33, 4).Value = "org.andante.enums.KafkaTopic"
$ws.Cells.Item(35, 2).Value = "`$assertionsDisabled"
$ws.Cells.Item(35, 4).Value = "boolean"
$ws.Cells.Item(36, 2).Value = "serialVersionUID"
$ws.Cells.Item(36, 4).Value = "long"
$ws.Cells.Item(37, 2).Value = "EQUAL"
$ws.Cells.Item(37, 3).Value = "public"
$ws.Cells.Item(37, 4).Value = "org.andante.rsql.operator.RSQLSearchOperator"
$ws.Cells.Item(39, 2).Value = "LESS_THAN_OR_EQUAL"
$ws.Cells.Item(40, 2).Value = "NOT_IN"
$ws.Cells.Item(41, 2).Value = "GREATER_THAN_OR_EQUAL"
$ws.Cells.Item(42, 2).Value = "GREATER_THAN"
$ws.Cells.Item(42, 3).Value = "public"
$ws.Cells.Item(42, 4).Value = "org.andante.rsql.operator.RSQLSearchOperator"
$ws.Cells.Item(43, 2).Value = "IN"
$ws.Cells.Item(44, 2).Value = "operator"
$ws.Cells.Item(44, 3).Value = "private"
$ws.Cells.Item(44, 4).Value = "cz.jirutka.rsql.parser.ast.ComparisonOperator"
$ws.Cells.Item(45, 2).Value = "LESS_THAN"
$ws.Cells.Item(46, 2).Value = "`$VALUES"
$ws.Cells.Item(46, 3).Value = "private"
$ws.Cells.Item(46, 4).Value = "org.andante.rsql.operator.RSQLSearchOperator[]"
$ws.Cells.Item(50, 2).Value = "MODIFY"
$ws.Cells.Item(50, 3).Value = "public"
$ws.Cells.Item(50, 4).Value = "org.andante.enums.OperationType"
$ws.Cells.Item(52, 2).Value = "`$VALUES"
$ws.Cells.Item(52, 3).Value = "private"
$ws.Cells.Item(52, 4).Value = "org.andante.enums.OperationType[]"
$ws.Cells.Item(53, 2).Value = "DATABASE_WILDCARD"
$ws.Cells.Item(53, 4).Value = "java.lang.String"
$ws.Cells.Item(54, 2).Value = "property"
$ws.Cells.Item(56, 2).Value = "RSQL_MAPPING_ERROR_MESSAGE"
$ws.Cells.Item(57, 2).Value = "operator"
$ws.Cells.Item(57, 4).Value = "cz.jirutka.rsql.parser.ast.ComparisonOperator"
$ws.Cells.Item(58, 2).Value = "RSQL_WILDCARD"
$ws.Cells.Item(60, 2).Value = "QUOTE"
$ws.Cells.Item(61, 2).Value = "LESS_THAN"
$ws.Cells.Item(61, 3).Value = "public"
$ws.Cells.Item(61, 4).Value = "org.andante.rsql.operator.RSQLReservedOperator"
$ws.Cells.Item(63, 2).Value = "SEMICOLON"
$ws.Cells.Item(64, 2).Value = "EXCLAMATION_MARK"
$ws.Cells.Item(65, 2).Value = "LEFT_BRACKET"
$ws.Cells.Item(66, 2).Value = "EQUALS"
$ws.Cells.Item(67, 2).Value = "APOSTROPHE"
$ws.Cells.Item(68, 2).Value = "COMMA"
$ws.Cells.Item(69, 2).Value = "RIGHT_BRACKET"
$ws.Cells.Item(70, 2).Value = "`$VALUES"
$ws.Cells.Item(70, 3).Value = "private"
$ws.Cells.Item(70, 4).Value = "org.andante.rsql.operator.RSQLReservedOperator[]"
$ws.Cells.Item(71, 2).Value = "reservedCharacter"
$ws.Cells.Item(71, 3).Value = "private"
$ws.Cells.Item(71, 4).Value = "java.lang.String"
$ws.Cells.Item(72, 2).Value = "TILDE"
$ws.Cells.Item(72, 3).Value = "public"
$ws.Cells.Item(72, 4).Value = "org.andante.rsql.operator.RSQLReservedOperator"
